$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers F1/G1, copying the header style (s="1") from E1 via format-only paste
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "TVN"
$ws.Range("G1").Value = "CTC"

# Rewrite data rows 2-51 (A:G) with the refreshed top-50 snippet data
$ws.Cells.Item(2,1).Value = 39
$ws.Cells.Item(2,2).Value = "5:45 AM"
$ws.Cells.Item(2,3).Value = 1169.48
$ws.Cells.Item(2,4).Value = 1199.48
$ws.Cells.Item(2,5).Value = 21.87
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 2
$ws.Cells.Item(3,1).Value = 49
$ws.Cells.Item(3,2).Value = "5:50 AM"
$ws.Cells.Item(3,3).Value = 1467.18
$ws.Cells.Item(3,4).Value = 1497.18
$ws.Cells.Item(3,5).Value = 24.93
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(4,1).Value = 91
$ws.Cells.Item(4,2).Value = "6:11 AM"
$ws.Cells.Item(4,3).Value = 2725.88
$ws.Cells.Item(4,4).Value = 2755.88
$ws.Cells.Item(4,5).Value = 22.16
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(5,1).Value = 100
$ws.Cells.Item(5,2).Value = "6:15 AM"
$ws.Cells.Item(5,3).Value = 2986.3175
$ws.Cells.Item(5,4).Value = 3016.3175
$ws.Cells.Item(5,5).Value = 30.9975
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 1.25
$ws.Cells.Item(6,1).Value = 120
$ws.Cells.Item(6,2).Value = "6:25 AM"
$ws.Cells.Item(6,3).Value = 3577.16
$ws.Cells.Item(6,4).Value = 3607.16
$ws.Cells.Item(6,5).Value = 49.37
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(7,1).Value = 151
$ws.Cells.Item(7,2).Value = "6:40 AM"
$ws.Cells.Item(7,3).Value = 4506.15
$ws.Cells.Item(7,4).Value = 4536.15
$ws.Cells.Item(7,5).Value = 33.81
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 1
$ws.Cells.Item(8,1).Value = 174
$ws.Cells.Item(8,2).Value = "6:52 AM"
$ws.Cells.Item(8,3).Value = 5193.04
$ws.Cells.Item(8,4).Value = 5223.04
$ws.Cells.Item(8,5).Value = 71.97
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(9,1).Value = 182
$ws.Cells.Item(9,2).Value = "6:56 AM"
$ws.Cells.Item(9,3).Value = 5440.89
$ws.Cells.Item(9,4).Value = 5470.89
$ws.Cells.Item(9,5).Value = 108.74
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(10,1).Value = 230
$ws.Cells.Item(10,2).Value = "7:20 AM"
$ws.Cells.Item(10,3).Value = 6880.786667
$ws.Cells.Item(10,4).Value = 6910.786667
$ws.Cells.Item(10,5).Value = 32.193333
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 1.333333
$ws.Cells.Item(11,1).Value = 239
$ws.Cells.Item(11,2).Value = "7:24 AM"
$ws.Cells.Item(11,3).Value = 7146.62
$ws.Cells.Item(11,4).Value = 7176.62
$ws.Cells.Item(11,5).Value = 32.95
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 2
$ws.Cells.Item(12,1).Value = 257
$ws.Cells.Item(12,2).Value = "7:33 AM"
$ws.Cells.Item(12,3).Value = 7691.65
$ws.Cells.Item(12,4).Value = 7721.65
$ws.Cells.Item(12,5).Value = 35.775
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(13,1).Value = 275
$ws.Cells.Item(13,2).Value = "7:43 AM"
$ws.Cells.Item(13,3).Value = 8248.53
$ws.Cells.Item(13,4).Value = 8278.53
$ws.Cells.Item(13,5).Value = 21.35
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(14,1).Value = 331
$ws.Cells.Item(14,2).Value = "8:10 AM"
$ws.Cells.Item(14,3).Value = 9903.9
$ws.Cells.Item(14,4).Value = 9933.9
$ws.Cells.Item(14,5).Value = 66.37
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(15,1).Value = 524
$ws.Cells.Item(15,2).Value = "9:47 AM"
$ws.Cells.Item(15,3).Value = 15702.116667
$ws.Cells.Item(15,4).Value = 15732.116667
$ws.Cells.Item(15,5).Value = 35.09
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(16,1).Value = 550
$ws.Cells.Item(16,2).Value = "10:00 A"
$ws.Cells.Item(16,3).Value = 16484
$ws.Cells.Item(16,4).Value = 16514
$ws.Cells.Item(16,5).Value = 26.08
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 1
$ws.Cells.Item(17,1).Value = 570
$ws.Cells.Item(17,2).Value = "10:10 A"
$ws.Cells.Item(17,3).Value = 17087.03
$ws.Cells.Item(17,4).Value = 17117.03
$ws.Cells.Item(17,5).Value = 37.05
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(18,1).Value = 993
$ws.Cells.Item(18,2).Value = "1:41 PM"
$ws.Cells.Item(18,3).Value = 29775.52
$ws.Cells.Item(18,4).Value = 29805.52
$ws.Cells.Item(18,5).Value = 41.35
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(19,1).Value = 1016
$ws.Cells.Item(19,2).Value = "1:53 PM"
$ws.Cells.Item(19,3).Value = 30454.73
$ws.Cells.Item(19,4).Value = 30484.73
$ws.Cells.Item(19,5).Value = 24.88
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 2
$ws.Cells.Item(20,1).Value = 1027
$ws.Cells.Item(20,2).Value = "1:59 PM"
$ws.Cells.Item(20,3).Value = 30807.17
$ws.Cells.Item(20,4).Value = 30837.17
$ws.Cells.Item(20,5).Value = 20.12
$ws.Cells.Item(20,6).Value = 0
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(21,1).Value = 1036
$ws.Cells.Item(21,2).Value = "2:03 PM"
$ws.Cells.Item(21,3).Value = 31058.49
$ws.Cells.Item(21,4).Value = 31088.49
$ws.Cells.Item(21,5).Value = 60.24
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 1
$ws.Cells.Item(22,1).Value = 1078
$ws.Cells.Item(22,2).Value = "2:24 PM"
$ws.Cells.Item(22,3).Value = 32329.13
$ws.Cells.Item(22,4).Value = 32359.13
$ws.Cells.Item(22,5).Value = 18.55
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 1
$ws.Cells.Item(23,1).Value = 1086
$ws.Cells.Item(23,2).Value = "2:28 PM"
$ws.Cells.Item(23,3).Value = 32575.77
$ws.Cells.Item(23,4).Value = 32605.77
$ws.Cells.Item(23,5).Value = 21.33
$ws.Cells.Item(23,6).Value = 0
$ws.Cells.Item(23,7).Value = 2
$ws.Cells.Item(24,1).Value = 1103
$ws.Cells.Item(24,2).Value = "2:36 PM"
$ws.Cells.Item(24,3).Value = 33070.88
$ws.Cells.Item(24,4).Value = 33100.88
$ws.Cells.Item(24,5).Value = 25.1
$ws.Cells.Item(24,6).Value = 0
$ws.Cells.Item(24,7).Value = 1
$ws.Cells.Item(25,1).Value = 1141
$ws.Cells.Item(25,2).Value = "2:55 PM"
$ws.Cells.Item(25,3).Value = 34219.87
$ws.Cells.Item(25,4).Value = 34249.87
$ws.Cells.Item(25,5).Value = 21.48
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 2
$ws.Cells.Item(26,1).Value = 1185
$ws.Cells.Item(26,2).Value = "3:17 PM"
$ws.Cells.Item(26,3).Value = 35530.52
$ws.Cells.Item(26,4).Value = 35560.52
$ws.Cells.Item(26,5).Value = 115.925
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 3
$ws.Cells.Item(27,1).Value = 1196
$ws.Cells.Item(27,2).Value = "3:23 PM"
$ws.Cells.Item(27,3).Value = 35858.715
$ws.Cells.Item(27,4).Value = 35888.715
$ws.Cells.Item(27,5).Value = 56.205
$ws.Cells.Item(27,6).Value = 0
$ws.Cells.Item(27,7).Value = 2
$ws.Cells.Item(28,1).Value = 1204
$ws.Cells.Item(28,2).Value = "3:27 PM"
$ws.Cells.Item(28,3).Value = 36097.916667
$ws.Cells.Item(28,4).Value = 36127.916667
$ws.Cells.Item(28,5).Value = 82.06
$ws.Cells.Item(28,6).Value = 0
$ws.Cells.Item(28,7).Value = 1.333333
$ws.Cells.Item(29,1).Value = 1214
$ws.Cells.Item(29,2).Value = "3:32 PM"
$ws.Cells.Item(29,3).Value = 36391.61
$ws.Cells.Item(29,4).Value = 36421.61
$ws.Cells.Item(29,5).Value = 96.52
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 7
$ws.Cells.Item(30,1).Value = 1221
$ws.Cells.Item(30,2).Value = "3:35 PM"
$ws.Cells.Item(30,3).Value = 36616.4
$ws.Cells.Item(30,4).Value = 36646.4
$ws.Cells.Item(30,5).Value = 23.47
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 0
$ws.Cells.Item(31,1).Value = 1245
$ws.Cells.Item(31,2).Value = "3:47 PM"
$ws.Cells.Item(31,3).Value = 37332.09
$ws.Cells.Item(31,4).Value = 37362.09
$ws.Cells.Item(31,5).Value = 193.12
$ws.Cells.Item(31,6).Value = 0
$ws.Cells.Item(31,7).Value = 6
$ws.Cells.Item(32,1).Value = 1257
$ws.Cells.Item(32,2).Value = "3:53 PM"
$ws.Cells.Item(32,3).Value = 37684.8
$ws.Cells.Item(32,4).Value = 37714.8
$ws.Cells.Item(32,5).Value = 36.82
$ws.Cells.Item(32,6).Value = 0
$ws.Cells.Item(32,7).Value = 1
$ws.Cells.Item(33,1).Value = 1273
$ws.Cells.Item(33,2).Value = "4:02 PM"
$ws.Cells.Item(33,3).Value = 38186.04
$ws.Cells.Item(33,4).Value = 38216.04
$ws.Cells.Item(33,5).Value = 28.77
$ws.Cells.Item(33,6).Value = 0
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(34,1).Value = 1283
$ws.Cells.Item(34,2).Value = "4:06 PM"
$ws.Cells.Item(34,3).Value = 38479.19
$ws.Cells.Item(34,4).Value = 38509.19
$ws.Cells.Item(34,5).Value = 51.24
$ws.Cells.Item(34,6).Value = 0
$ws.Cells.Item(34,7).Value = 3
$ws.Cells.Item(35,1).Value = 1291
$ws.Cells.Item(35,2).Value = "4:11 PM"
$ws.Cells.Item(35,3).Value = 38728.38
$ws.Cells.Item(35,4).Value = 38758.38
$ws.Cells.Item(35,5).Value = 100.57
$ws.Cells.Item(35,6).Value = 0
$ws.Cells.Item(35,7).Value = 2
$ws.Cells.Item(36,1).Value = 1302
$ws.Cells.Item(36,2).Value = "4:16 PM"
$ws.Cells.Item(36,3).Value = 39032.49
$ws.Cells.Item(36,4).Value = 39062.49
$ws.Cells.Item(36,5).Value = 230.5
$ws.Cells.Item(36,6).Value = 0
$ws.Cells.Item(36,7).Value = 1
$ws.Cells.Item(37,1).Value = 1311
$ws.Cells.Item(37,2).Value = "4:20 PM"
$ws.Cells.Item(37,3).Value = 39315.5
$ws.Cells.Item(37,4).Value = 39345.5
$ws.Cells.Item(37,5).Value = 76.02
$ws.Cells.Item(37,6).Value = 0
$ws.Cells.Item(37,7).Value = 1.5
$ws.Cells.Item(38,1).Value = 1318
$ws.Cells.Item(38,2).Value = "4:24 PM"
$ws.Cells.Item(38,3).Value = 39526.18
$ws.Cells.Item(38,4).Value = 39556.18
$ws.Cells.Item(38,5).Value = 111.31
$ws.Cells.Item(38,6).Value = 0
$ws.Cells.Item(38,7).Value = 2.5
$ws.Cells.Item(39,1).Value = 1335
$ws.Cells.Item(39,2).Value = "4:32 PM"
$ws.Cells.Item(39,3).Value = 40027.72
$ws.Cells.Item(39,4).Value = 40057.72
$ws.Cells.Item(39,5).Value = 47.705
$ws.Cells.Item(39,6).Value = 0
$ws.Cells.Item(39,7).Value = 0.5
$ws.Cells.Item(40,1).Value = 1345
$ws.Cells.Item(40,2).Value = "4:38 PM"
$ws.Cells.Item(40,3).Value = 40347.76
$ws.Cells.Item(40,4).Value = 40377.76
$ws.Cells.Item(40,5).Value = 94.74
$ws.Cells.Item(40,6).Value = 0
$ws.Cells.Item(40,7).Value = 1
$ws.Cells.Item(41,1).Value = 1357
$ws.Cells.Item(41,2).Value = "4:43 PM"
$ws.Cells.Item(41,3).Value = 40690.25
$ws.Cells.Item(41,4).Value = 40720.25
$ws.Cells.Item(41,5).Value = 163.16
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 5
$ws.Cells.Item(42,1).Value = 1372
$ws.Cells.Item(42,2).Value = "4:51 PM"
$ws.Cells.Item(42,3).Value = 41141.465
$ws.Cells.Item(42,4).Value = 41171.465
$ws.Cells.Item(42,5).Value = 98.08
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 5
$ws.Cells.Item(43,1).Value = 1391
$ws.Cells.Item(43,2).Value = "5:01 PM"
$ws.Cells.Item(43,3).Value = 41726.28
$ws.Cells.Item(43,4).Value = 41756.28
$ws.Cells.Item(43,5).Value = 40.32
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 2
$ws.Cells.Item(44,1).Value = 1412
$ws.Cells.Item(44,2).Value = "5:11 PM"
$ws.Cells.Item(44,3).Value = 42334.07
$ws.Cells.Item(44,4).Value = 42364.07
$ws.Cells.Item(44,5).Value = 101.66
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 2
$ws.Cells.Item(45,1).Value = 1427
$ws.Cells.Item(45,2).Value = "5:19 PM"
$ws.Cells.Item(45,3).Value = 42809.97
$ws.Cells.Item(45,4).Value = 42839.97
$ws.Cells.Item(45,5).Value = 161.82
$ws.Cells.Item(45,6).Value = 0
$ws.Cells.Item(45,7).Value = 9
$ws.Cells.Item(46,1).Value = 1441
$ws.Cells.Item(46,2).Value = "5:26 PM"
$ws.Cells.Item(46,3).Value = 43229.93
$ws.Cells.Item(46,4).Value = 43259.93
$ws.Cells.Item(46,5).Value = 24.11
$ws.Cells.Item(46,6).Value = 0
$ws.Cells.Item(46,7).Value = 1
$ws.Cells.Item(47,1).Value = 1451
$ws.Cells.Item(47,2).Value = "5:30 PM"
$ws.Cells.Item(47,3).Value = 43516.02
$ws.Cells.Item(47,4).Value = 43546.02
$ws.Cells.Item(47,5).Value = 88.58
$ws.Cells.Item(47,6).Value = 0
$ws.Cells.Item(47,7).Value = 2
$ws.Cells.Item(48,1).Value = 1479
$ws.Cells.Item(48,2).Value = "5:44 PM"
$ws.Cells.Item(48,3).Value = 44340.68
$ws.Cells.Item(48,4).Value = 44370.68
$ws.Cells.Item(48,5).Value = 49.65
$ws.Cells.Item(48,6).Value = 0
$ws.Cells.Item(48,7).Value = 3
$ws.Cells.Item(49,1).Value = 1492
$ws.Cells.Item(49,2).Value = "5:51 PM"
$ws.Cells.Item(49,3).Value = 44738.87
$ws.Cells.Item(49,4).Value = 44768.87
$ws.Cells.Item(49,5).Value = 22.62
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(50,1).Value = 1522
$ws.Cells.Item(50,2).Value = "6:06 PM"
$ws.Cells.Item(50,3).Value = 45637.495
$ws.Cells.Item(50,4).Value = 45667.495
$ws.Cells.Item(50,5).Value = 33.28
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 0.5
$ws.Cells.Item(51,1).Value = 1551
$ws.Cells.Item(51,2).Value = "6:20 PM"
$ws.Cells.Item(51,3).Value = 46502.28
$ws.Cells.Item(51,4).Value = 46532.28
$ws.Cells.Item(51,5).Value = 39.27
$ws.Cells.Item(51,6).Value = 1.01
$ws.Cells.Item(51,7).Value = 1
